# Update the "想去人数" (interest count) figures in the 展览 and 全部类型
# sheets to reflect a fresh scrape, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 40
$ws1.Range("F7").Value  = 14509
$ws1.Range("F9").Value  = 660
$ws1.Range("F10").Value = 15124
$ws1.Range("F12").Value = 8543
$ws1.Range("F29").Value = 412
$ws1.Range("F33").Value = 254
$ws1.Range("F36").Value = 5257

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 40
$ws4.Range("F7").Value  = 14509
$ws4.Range("F9").Value  = 660
$ws4.Range("F10").Value = 15124
$ws4.Range("F12").Value = 8543
$ws4.Range("F32").Value = 412
$ws4.Range("F36").Value = 254
$ws4.Range("F39").Value = 5257
